# Generate Report for Handback
# Updates the handback timestamp values that were refreshed when the
# localization report was regenerated.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for the
# 5cb1b9c8-...md row moves forward in time.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-07 17:08:56"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the 5cb1b9c8-... row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-07 17:08:50"
$wsZhCn.Range("K3").Value = "2016-09-07 17:09:54"

# "de-de" sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the 5cb1b9c8-... row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-07 17:08:56"
$wsDeDe.Range("K3").Value = "2016-09-07 17:10:08"
